$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.042.87"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "'1.662.13"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'216.32"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.252"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D10").Value = "'20.16"
$ws.Range("E10").Value = "  +4.85%  "
$ws.Range("D11").Value = "'0.0884"
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("D12").Value = "'1.894.17"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "'1.665.42"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'4.10"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").Value = "'65.73"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "'27.068.68"
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "'236.83"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'0.0₃0739"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.30"
$ws.Range("E23").Value = "  +2.62%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.23"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "'145.46"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'0.0499"
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("D32").Value = "'1.553.87"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").Value = "  +4.41%  "
$ws.Range("E35").Value = "  +6.53%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.578"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'0.907"
$ws.Range("E38").Value = "  +9.36%  "
$ws.Range("D39").Value = "'0.0170"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").Value = "'6.08"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'66.78"
$ws.Range("E42").Value = "  +8.70%  "
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("D44").Value = "'0.971"
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("D45").Value = "'1.804.46"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "'0.779"
$ws.Range("E46").Value = "  +2.37%  "
$ws.Range("D47").Value = "'90.65"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("D48").Value = "'1.54"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  +2.55%  "
